$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on D-column cells whose new values would
# otherwise be auto-parsed as numbers by Excel (losing exact text/trailing zeros).
$textFormatRows = @(5,7,8,9,10,12,13,14,15,16,20,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $textFormatRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '25.835.52'
$ws.Range("E2").Value = '  +0.17%  '

# Row 3
$ws.Range("D3").Value = '1.735.99'
$ws.Range("E3").Value = '  -0.66%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '230.68'
$ws.Range("E5").Value = '  -2.10%  '

# Row 6
$ws.Range("E6").Value = '  +0.03%  '

# Row 7
$ws.Range("D7").Value = '0.5138'
$ws.Range("E7").Value = '  +0.93%  '

# Row 8
$ws.Range("D8").Value = '0.2794'
$ws.Range("E8").Value = '  +4.32%  '

# Row 9
$ws.Range("D9").Value = '39.31'
$ws.Range("E9").Value = '  -2.50%  '

# Row 10
$ws.Range("D10").Value = '0.06100'
$ws.Range("E10").Value = '  -1.16%  '

# Row 11
$ws.Range("D11").Value = '1.754.16'
$ws.Range("E11").Value = '  +0.12%  '

# Row 12
$ws.Range("D12").Value = '0.07045'
$ws.Range("E12").Value = '  +1.54%  '

# Row 13
$ws.Range("D13").Value = '15.24'
$ws.Range("E13").Value = '  -0.61%  '

# Row 14
$ws.Range("D14").Value = '0.6430'
$ws.Range("E14").Value = '  +2.90%  '

# Row 15
$ws.Range("D15").Value = '4.514'
$ws.Range("E15").Value = '  +0.93%  '

# Row 16
$ws.Range("D16").Value = '76.73'
$ws.Range("E16").Value = '  -1.46%  '

# Row 17
$ws.Range("E17").Value = '  +0.13%  '

# Row 18
$ws.Range("E18").Value = '  -0.05%  '

# Row 19
$ws.Range("D19").Value = '25.828.37'
$ws.Range("E19").Value = '  +0.06%  '

# Row 20
$ws.Range("D20").Value = '11.47'
$ws.Range("E20").Value = '  -0.92%  '

# Row 21
$ws.Range("E21").Value = '  -0.47%  '

# Row 22
$ws.Range("D22").Value = '1.973.27'
$ws.Range("E22").Value = '  +0.66%  '

# Row 23
$ws.Range("D23").Value = '4.138'
$ws.Range("E23").Value = '  +2.43%  '

# Row 24
$ws.Range("D24").Value = '8.661'
$ws.Range("E24").Value = '  +4.94%  '

# Row 25
$ws.Range("D25").Value = '5.131'
$ws.Range("E25").Value = '  -0.03%  '

# Row 26
$ws.Range("D26").Value = '139.37'
$ws.Range("E26").Value = '  +2.23%  '

# Row 27
$ws.Range("D27").Value = '1.514'
$ws.Range("E27").Value = '  +3.38%  '

# Row 28
$ws.Range("D28").Value = '15.04'
$ws.Range("E28").Value = '  -0.12%  '

# Row 29
$ws.Range("D29").Value = '1.798'
$ws.Range("E29").Value = '  +1.81%  '

# Row 30
$ws.Range("D30").Value = '102.13'
$ws.Range("E30").Value = '  -0.27%  '

# Row 31
$ws.Range("D31").Value = '0.08322'
$ws.Range("E31").Value = '  +1.78%  '

# Row 32
$ws.Range("D32").Value = '3.673'
$ws.Range("E32").Value = '  -0.03%  '

# Row 33
$ws.Range("D33").Value = '3.416'
$ws.Range("E33").Value = '  +0.79%  '

# Row 34
$ws.Range("D34").Value = '0.04489'
$ws.Range("E34").Value = '  +2.57%  '

# Row 35
$ws.Range("D35").Value = '2.613'
$ws.Range("E35").Value = '  -1.34%  '

# Row 36
$ws.Range("D36").Value = '0.9825'
$ws.Range("E36").Value = '  -1.38%  '

# Row 37
$ws.Range("D37").Value = '0.6104'
$ws.Range("E37").Value = '  +2.09%  '

# Row 38
$ws.Range("D38").Value = '2.641'
$ws.Range("E38").Value = '  +1.86%  '

# Row 39
$ws.Range("D39").Value = '0.01578'
$ws.Range("E39").Value = '  +1.36%  '

# Row 40
$ws.Range("D40").Value = '1.935'
$ws.Range("E40").Value = '  +0.59%  '

# Row 41
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  +0.00%  '

# Row 42
$ws.Range("D42").Value = '100.37'
$ws.Range("E42").Value = '  -1.25%  '

# Row 43
$ws.Range("D43").Value = '0.3828'
$ws.Range("E43").Value = '  +0.01%  '

# Row 44
$ws.Range("D44").Value = '0.7252'
$ws.Range("E44").Value = '  -2.70%  '

# Row 45
$ws.Range("D45").Value = '4.969'
$ws.Range("E45").Value = '  +1.68%  '

# Row 46
$ws.Range("D46").Value = '0.05398'
$ws.Range("E46").Value = '  -1.70%  '

# Row 47
$ws.Range("D47").Value = '6.267'
$ws.Range("E47").Value = '  +5.78%  '

# Row 48
$ws.Range("D48").Value = '0.1115'
$ws.Range("E48").Value = '  +1.84%  '

# Row 49
$ws.Range("D49").Value = '52.93'
$ws.Range("E49").Value = '  +0.98%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '29.94'
$ws.Range("E50").Value = '  -0.30%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.637'
$ws.Range("E51").Value = '  +3.08%  '
